$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 411
$ws1.Range("F5").Value = 5471
$ws1.Range("F6").Value = 5471
$ws1.Range("F7").Value = 202
$ws1.Range("F9").Value = 547
$ws1.Range("F11").Value = 1229
$ws1.Range("F14").Value = 797
$ws1.Range("F15").Value = 6477
$ws1.Range("F18").Value = 117
$ws1.Range("F19").Value = 4485
$ws1.Range("F22").Value = 4149
$ws1.Range("F23").Value = 4073
$ws1.Range("F24").Value = 200
$ws1.Range("F25").Value = 204
$ws1.Range("F26").Value = 273
$ws1.Range("F27").Value = 262
$ws1.Range("F34").Value = 7336
$ws1.Range("F36").Value = 1219
$ws1.Range("F37").Value = 605
$ws1.Range("F41").Value = 1467
$ws1.Range("F43").Value = 807
$ws1.Range("F45").Value = 3534
$ws1.Range("F47").Value = 9
$ws1.Range("F49").Value = 809
$ws1.Range("F50").Value = 1013

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 4
$ws2.Range("F13").Value = 17
$ws2.Range("F14").Value = 148
$ws2.Range("F19").Value = 63
$ws2.Range("F22").Value = 847

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 411
$ws4.Range("F5").Value = 5471
$ws4.Range("F6").Value = 5471
$ws4.Range("F7").Value = 202
$ws4.Range("F10").Value = 547
$ws4.Range("F12").Value = 1229
$ws4.Range("F15").Value = 4
$ws4.Range("F16").Value = 797
$ws4.Range("F17").Value = 6477
$ws4.Range("F20").Value = 117
$ws4.Range("F21").Value = 4485
$ws4.Range("F24").Value = 4149
$ws4.Range("F25").Value = 4073
$ws4.Range("F26").Value = 200
$ws4.Range("F27").Value = 204
$ws4.Range("F28").Value = 273
$ws4.Range("F29").Value = 262
$ws4.Range("F33").Value = 148
$ws4.Range("F34").Value = 7336
$ws4.Range("F36").Value = 1219
$ws4.Range("F37").Value = 605
$ws4.Range("F41").Value = 1467
$ws4.Range("F43").Value = 807
$ws4.Range("F45").Value = 3534
$ws4.Range("F48").Value = 809
$ws4.Range("F49").Value = 1013
